$wb     = $excel.ActiveWorkbook
$matrix = $wb.Worksheets.Item("matrix")
$long   = $wb.Worksheets.Item("long")

# ---------------------------------------------------------------------
# 1. Rename the "Final uses" / "FU" final-demand column to "Final
#    consumption expenditure of households" / "FCH" everywhere it is
#    used: the matrix sheet's row-6/row-7 header cells, and every row
#    of the long sheet's lookup column (D) that pointed at the old
#    "FU" code.
# ---------------------------------------------------------------------
$matrix.Range("E6:F6").Value = "Final consumption expenditure of households"
$matrix.Range("E7:F7").Value = "FCH"
$long.Range("D2:D17").Value  = "FCH"

# ---------------------------------------------------------------------
# 2. matrix!H8:H15 ("<country><product>" key column) becomes a single
#    shared formula group - writing the formula across the whole
#    range in one assignment lets the relative references adjust per
#    row while keeping the formula itself shared.
# ---------------------------------------------------------------------
$matrix.Range("H8:H15").Formula = "=B8&D8"

# ---------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping: "matrix" becomes the active
#    tab (selection resting on F7); "long" is no longer active
#    (selection resting on D4).
# ---------------------------------------------------------------------
[void]$long.Range("D4").Select()
[void]$matrix.Activate()
[void]$matrix.Range("F7").Select()
